$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.005.95'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '3.257.75'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.30'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.599'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('E9').Value = '  -3.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.59'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('D12').Value = '3.826.36'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('E13').Value = '  +1.37%  '
$ws.Range('D14').Value = '68.008.26'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.21'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.21%  '
$ws.Range('E16').Value = '  -2.97%  '
$ws.Range('D17').Value = '3.268.47'
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('E18').Value = '  -3.21%  '
$ws.Range('E19').Value = '  -3.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '414.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.93%  '
$ws.Range('E21').Value = '  -3.18%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.94'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.95%  '
$ws.Range('E24').Value = '  -2.94%  '
$ws.Range('E25').Value = '  -4.18%  '
$ws.Range('E26').Value = '  -1.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.26'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.996'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.46%  '
$ws.Range('E29').Value = '  -2.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.55'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.26%  '
$ws.Range('E31').Value = '  -6.17%  '
$ws.Range('E32').Value = '  -5.13%  '
$ws.Range('E33').Value = '  -5.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '164.59'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('E35').Value = '  -6.07%  '
$ws.Range('E36').Value = '  -5.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.47'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.790'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.37%  '
$ws.Range('E39').Value = '  -4.35%  '
$ws.Range('E40').Value = '  -5.00%  '
$ws.Range('D41').Value = '2.618.27'
$ws.Range('E41').Value = '  -1.60%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0671'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.19%  '
$ws.Range('E43').Value = '  -5.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '333.06'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '24.12'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.02%  '
$ws.Range('E46').Value = '  -4.39%  '
$ws.Range('E48').Value = '  -2.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0997'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.40%  '
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '30.39'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.80%  '
